$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column stays text (matches source data which stores formatted price strings)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '30.161.80'
$ws.Cells.Item(2, 5).Value = '  +1.06%  '
$ws.Cells.Item(3, 4).Value = '1.894.43'
$ws.Cells.Item(3, 5).Value = '  +0.36%  '
$ws.Cells.Item(4, 4).Value = '1.003'
$ws.Cells.Item(4, 5).Value = '  +0.26%  '
$ws.Cells.Item(5, 4).Value = '0.7406'
$ws.Cells.Item(5, 5).Value = '  -0.81%  '
$ws.Cells.Item(6, 4).Value = '243.17'
$ws.Cells.Item(6, 5).Value = '  +0.32%  '
$ws.Cells.Item(7, 4).Value = '1.002'
$ws.Cells.Item(7, 5).Value = '  +0.15%  '
$ws.Cells.Item(8, 4).Value = '0.3173'
$ws.Cells.Item(8, 5).Value = '  +1.84%  '
$ws.Cells.Item(9, 4).Value = '0.07226'
$ws.Cells.Item(9, 5).Value = '  +1.33%  '
$ws.Cells.Item(10, 4).Value = '24.98'
$ws.Cells.Item(10, 5).Value = '  -1.13%  '
$ws.Cells.Item(11, 4).Value = '0.08358'
$ws.Cells.Item(11, 5).Value = '  -1.61%  '
$ws.Cells.Item(12, 2).Value = 'WrappedEther'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(12, 4).Value = '1.955.42'
$ws.Cells.Item(12, 5).Value = '  +3.43%  '
$ws.Cells.Item(13, 2).Value = 'Polygon'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(13, 4).Value = '0.7618'
$ws.Cells.Item(13, 5).Value = '  +0.31%  '
$ws.Cells.Item(14, 2).Value = 'Polkadot'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(14, 4).Value = '5.461'
$ws.Cells.Item(14, 5).Value = '  +2.04%  '
$ws.Cells.Item(15, 4).Value = '93.06'
$ws.Cells.Item(15, 5).Value = '  -0.28%  '
$ws.Cells.Item(16, 4).Value = '6.171'
$ws.Cells.Item(16, 5).Value = '  +0.01%  '
$ws.Cells.Item(17, 4).Value = '30.190.85'
$ws.Cells.Item(17, 5).Value = '  +1.14%  '
$ws.Cells.Item(18, 4).Value = '250.98'
$ws.Cells.Item(18, 5).Value = '  +3.15%  '
$ws.Cells.Item(19, 4).Value = '13.66'
$ws.Cells.Item(19, 5).Value = '  -0.29%  '
$ws.Cells.Item(20, 4).Value = '0.000007896'
$ws.Cells.Item(20, 5).Value = '  +1.28%  '
$ws.Cells.Item(21, 4).Value = '2.178.98'
$ws.Cells.Item(21, 5).Value = '  +1.14%  '
$ws.Cells.Item(22, 5).Value = '  +0.23%  '
$ws.Cells.Item(23, 4).Value = '7.970'
$ws.Cells.Item(23, 5).Value = '  -0.30%  '
$ws.Cells.Item(24, 4).Value = '1.003'
$ws.Cells.Item(25, 4).Value = '0.1585'
$ws.Cells.Item(25, 5).Value = '  -0.15%  '
$ws.Cells.Item(26, 4).Value = '9.319'
$ws.Cells.Item(26, 5).Value = '  -0.36%  '
$ws.Cells.Item(27, 4).Value = '164.86'
$ws.Cells.Item(27, 5).Value = '  +1.51%  '
$ws.Cells.Item(28, 4).Value = '18.80'
$ws.Cells.Item(28, 5).Value = '  +0.32%  '
$ws.Cells.Item(29, 4).Value = '2.068'
$ws.Cells.Item(29, 5).Value = '  +2.17%  '
$ws.Cells.Item(30, 4).Value = '1.484'
$ws.Cells.Item(30, 5).Value = '  -0.59%  '
$ws.Cells.Item(31, 4).Value = '4.591'
$ws.Cells.Item(31, 5).Value = '  +2.25%  '
$ws.Cells.Item(32, 4).Value = '1.535'
$ws.Cells.Item(32, 5).Value = '  +0.17%  '
$ws.Cells.Item(33, 4).Value = '4.214'
$ws.Cells.Item(33, 5).Value = '  +2.55%  '
$ws.Cells.Item(34, 4).Value = '0.05382'
$ws.Cells.Item(34, 5).Value = '  -0.48%  '
$ws.Cells.Item(35, 4).Value = '1.257'
$ws.Cells.Item(35, 5).Value = '  +1.54%  '
$ws.Cells.Item(36, 4).Value = '0.7891'
$ws.Cells.Item(36, 5).Value = '  +5.92%  '
$ws.Cells.Item(37, 4).Value = '1.004'
$ws.Cells.Item(37, 5).Value = '  +0.18%  '
$ws.Cells.Item(38, 4).Value = '2.732'
$ws.Cells.Item(38, 5).Value = '  +0.78%  '
$ws.Cells.Item(39, 5).Value = '  +1.97%  '
$ws.Cells.Item(40, 4).Value = '2.769'
$ws.Cells.Item(40, 5).Value = '  -0.10%  '
$ws.Cells.Item(41, 4).Value = '0.4576'
$ws.Cells.Item(41, 5).Value = '  +2.86%  '
$ws.Cells.Item(42, 4).Value = '1.100.69'
$ws.Cells.Item(42, 5).Value = '  +0.89%  '
$ws.Cells.Item(43, 4).Value = '6.078'
$ws.Cells.Item(43, 5).Value = '  +0.06%  '
$ws.Cells.Item(44, 4).Value = '72.84'
$ws.Cells.Item(44, 5).Value = '  +0.61%  '
$ws.Cells.Item(45, 4).Value = '0.8735'
$ws.Cells.Item(45, 5).Value = '  +2.23%  '
$ws.Cells.Item(46, 2).Value = 'Quant'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(46, 4).Value = '104.55'
$ws.Cells.Item(46, 5).Value = '  +2.24%  '
$ws.Cells.Item(47, 2).Value = 'PaxDollar'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(47, 4).Value = '1.003'
$ws.Cells.Item(47, 5).Value = '  +0.19%  '
$ws.Cells.Item(48, 4).Value = '1.868'
$ws.Cells.Item(48, 5).Value = '  +0.45%  '
$ws.Cells.Item(49, 4).Value = '7.615'
$ws.Cells.Item(49, 5).Value = '  -0.89%  '
$ws.Cells.Item(50, 4).Value = '9.644'
$ws.Cells.Item(50, 5).Value = '  -0.86%  '
$ws.Cells.Item(51, 4).Value = '2.067.75'
$ws.Cells.Item(51, 5).Value = '  +0.91%  '
